$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 78

$ws.Cells.Item($row, 1).Value = "CompaNanny"
$ws.Cells.Item($row, 2).Value = "CompaNanny Statenkwartier BSO"
$ws.Cells.Item($row, 3).Value = "VGO"

$dateCell = $ws.Cells.Item($row, 4)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2024-09-16"
$dateCell.Style = "Normal"

$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 1
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
